$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: Copy style (s="1") template from C65 to the C column of new rows 158-176 ---
$ws.Range("C65").Copy() | Out-Null
$ws.Range("C158:C176").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Step 2: Update existing cells whose text content changed ---
$ws.Range("D147").Value = "setItemMoney;buyItemConfirmedDialogData;buyItemConfirmedDialog"

$ws.Range("D148").Value = "buyItemWindow"

$ws.Range("A149").Value = "setItemMoney"

$ws.Range("D153").Value = "buyItemWindow"

$ws.Range("D155").Value = "costMoney;buyItemGainItem;cityHaveItem"

$ws.Range("A158").Value = "sellItem"
$ws.Range("B158").Value = "卖出道具"
$ws.Range("C158").Value = "eventList"
$ws.Range("D158").Value = "close;sellItemHaveItem"

$ws.Range("A159").Value = "sellItemHaveItem"
$ws.Range("C159").Value = "condition"
$ws.Range("D159").Value = "haveItem;sellItemStart;sellItemNoItem"

$ws.Range("A160").Value = "sellItemStart"
$ws.Range("C160").Value = "eventList"
$ws.Range("D160").Value = "sellItemWindow"

$ws.Range("A161").Value = "sellItemNoItem"
$ws.Range("C161").Value = "eventList"
$ws.Range("D161").Value = "sellItemNoItemDialog;shop"

$ws.Range("A162").Value = "sellItemNoItemDialog"
$ws.Range("C162").Value = "dialog"
$ws.Range("D162").Value = "dialog_i_dont_have_item"

$ws.Range("A163").Value = "sellItemWindow"
$ws.Range("C163").Value = "window"
$ws.Range("D163").Value = "ItemBrowsePanel;sellItemSelected;buyItemCanceled;1"

$ws.Range("A164").Value = "sellItemSelected"
$ws.Range("C164").Value = "eventList"
$ws.Range("D164").Value = "sellItemSmallWindow"

$ws.Range("A165").Value = "sellItemSmallWindow"
$ws.Range("C165").Value = "window"
$ws.Range("D165").Value = "ItemInfoPanel;sellItemConfirmed;sellItemCancelSmallWindow;1;reserved.itemId"

$ws.Range("A166").Value = "sellItemConfirmed"
$ws.Range("C166").Value = "eventList"
$ws.Range("D166").Value = "setItemMoney;buyItemConfirmedDialogData;selltemConfirmedDialog"

$ws.Range("A167").Value = "sellItemCancelSmallWindow"
$ws.Range("C167").Value = "eventList"
$ws.Range("D167").Value = "sellItemWindow"

$ws.Range("A168").Value = "selltemConfirmedDialog"
$ws.Range("C168").Value = "dialogYesNo"
$ws.Range("D168").Value = "dialog_sell_item;selltemDealStart;sellItemDealCanceled"

$ws.Range("A169").Value = "selltemDealStart"
$ws.Range("C169").Value = "eventList"
$ws.Range("D169").Value = "gainMoney;sellItemLoseItem;sellItemHaveItem"

$ws.Range("A170").Value = "sellItemDealCanceled"
$ws.Range("C170").Value = "eventList"
$ws.Range("D170").Value = "sellItemWindow"

$ws.Range("A171").Value = "gainMoney"
$ws.Range("C171").Value = "dataChange"
$ws.Range("D171").Value = "money;+;cache.money"

$ws.Range("A172").Value = "sellItemLoseItem"
$ws.Range("C172").Value = "dataChange"
$ws.Range("D172").Value = "item;sell;reserved.itemId"

$ws.Range("A173").Value = "arrangeSailor"
$ws.Range("B173").Value = "分配水手"
$ws.Range("C173").Value = "eventList"
$ws.Range("D173").Value = "arrangeSailorHaveShip"

$ws.Range("A174").Value = "arrangeSailorHaveShip"
$ws.Range("C174").Value = "condition"
$ws.Range("D174").Value = "haveShip;arrangeSailorWindow1;arrangeSailorFailedDialog"

$ws.Range("A175").Value = "arrangeSailorWindow1"
$ws.Range("C175").Value = "window"
$ws.Range("D175").Value = "SailorNumberPanel;;;0"

$ws.Range("A176").Value = "arrangeSailorFailedDialog"
$ws.Range("C176").Value = "dialog"
$ws.Range("D176").Value = "dialog_no_ship_no_game"
